# Op_Text.xlsx cleanup pass:
#  - add a dynamic-array LOWER() mirror of the existing A2:A224 pattern list
#    (spills into A225:A447)
#  - append lower-cased "replace" values in column B for the rows whose
#    upper/mixed-case counterparts already carry a replacement, plus a
#    couple of brand new find/replace pairs
#  - append one brand-new pair (row 448)
#  - leave the cursor/selection on B376, matching the author's last position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. dynamic array formula: lower-case mirror of the original list ---
$ws.Range("A225:A447").FormulaArray = "=LOWER(A2:A224)"

# --- 2. new lower-cased replacement values in column B ------------------
# (row -> old row it mirrors): 290<-67, 373<-150, 374<-151, 375<-152,
# 376<-153, 386<-163, 388<-165, 389<-166, 390<-167, 413<-190, 414<-191,
# 415<-192, 420<-197, 423<-200, 424<-201, 432<-209, 433<-210, 435<-212
$ws.Range("B290").Value = "trụ điện"
$ws.Range("B373").Value = "nguyễn văn đẹp"
$ws.Range("B374").Value = "tân phước"
$ws.Range("B375").Value = 9
$ws.Range("B376").Value = "d8"
$ws.Range("B386").Value = 883
$ws.Range("B388").Value = 1
$ws.Range("B389").Value = "17/1"
$ws.Range("B390").Value = "30/4"
$ws.Range("B413").Value = "khu vực"
$ws.Range("B414").Value = "khu vực"
$ws.Range("B415").Value = "khu vực"
$ws.Range("B420").Value = "tân khánh"
$ws.Range("B423").Value = "mai thị non"
$ws.Range("B424").Value = "9 &  5"
$ws.Range("B432").Value = "so 8"
$ws.Range("B433").Value = "số 2"
$ws.Range("B435").Value = "19/5"

# --- 3. brand new find/replace pair appended after the spill range ------
$ws.Range("A448").Value = "26 tháng 3"
$ws.Range("B448").Value = "26/3"

# --- 4. restore the author's last selection ------------------------------
$ws.Range("B376").Select()
